# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '58.023.34'
$ws.Range('E2').Value = '  -1.73%  '
$ws.Range('D3').Value = '2.461.85'
$ws.Range('E3').Value = '  -1.43%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '518.75'
$ws.Range('E5').Value = '  -2.99%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '132.89'
$ws.Range('E6').Value = '  -2.44%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.556'
$ws.Range('E8').Value = '  -1.75%  '
$ws.Range('D9').Value = '2.473.77'
$ws.Range('E9').Value = '  -1.32%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0975'
$ws.Range('E10').Value = '  -4.20%  '
$ws.Range('E11').Value = '  -0.83%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '5.30'
$ws.Range('E12').Value = '  -1.92%  '
$ws.Range('E13').Value = '  -3.35%  '
$ws.Range('D14').Value = '2.900.79'
$ws.Range('E14').Value = '  -1.34%  '
$ws.Range('D15').Value = '57.946.04'
$ws.Range('E15').Value = '  -1.66%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '21.88'
$ws.Range('E16').Value = '  -3.44%  '
$ws.Range('E17').Value = '  -3.05%  '
$ws.Range('D18').Value = '2.480.34'
$ws.Range('E18').Value = '  -1.05%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '10.59'
$ws.Range('E19').Value = '  -4.03%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '319.26'
$ws.Range('E20').Value = '  -1.33%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '4.13'
$ws.Range('E21').Value = '  -2.84%  '
$ws.Range('E22').Value = '  +0.03%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.69'
$ws.Range('E23').Value = '  -4.96%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '64.40'
$ws.Range('E24').Value = '  -1.22%  '
$ws.Range('E25').Value = '  -3.53%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.997'
$ws.Range('E26').Value = '  -0.12%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.160'
$ws.Range('E27').Value = '  -2.78%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '7.28'
$ws.Range('E28').Value = '  -3.06%  '
$ws.Range('D29').Value = '0.0₃0740'
$ws.Range('E29').Value = '  -2.82%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '168.39'
$ws.Range('E30').Value = '  -1.71%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.69'
$ws.Range('E31').Value = '  -3.44%  '
$ws.Range('E32').Value = '  -4.33%  '
$ws.Range('E33').Value = '  -1.08%  '
$ws.Range('E34').Value = '  -0.03%  '
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '17.99'
$ws.Range('E36').Value = '  -2.02%  '
$ws.Range('E37').Value = '  -2.40%  '
$ws.Range('E38').Value = '  -2.64%  '
$ws.Range('E39').Value = '  -1.03%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.46'
$ws.Range('E40').Value = '  -4.68%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.785'
$ws.Range('E41').Value = '  -1.66%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.42'
$ws.Range('E42').Value = '  -4.33%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '270.24'
$ws.Range('E43').Value = '  -4.41%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '4.96'
$ws.Range('E44').Value = '  -5.74%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.589'
$ws.Range('E45').Value = '  -2.22%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '123.16'
$ws.Range('E46').Value = '  -5.40%  '
$ws.Range('E47').Value = '  -1.90%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0484'
$ws.Range('E48').Value = '  -3.13%  '
$ws.Range('E49').Value = '  -3.07%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '16.76'
$ws.Range('E50').Value = '  -2.93%  '
$ws.Range('D51').Value = '1.731.91'
$ws.Range('E51').Value = '  -1.35%  '
